# The post about "「太陽」شَمْس" (originally worksheet row 5) was removed
# from the source data. Delete that entire row; Excel will automatically
# shift all subsequent rows up by one and update the sheet's used range /
# dimension accordingly (from A1:C209 to A1:C208).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Delete()
